$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Prepend "${1} ${2} ${3} ${4} ${5} ${6} " as its own run, right before the
#    existing "${7" run (paragraph that renders as "${7} ${8} ... ${20} ").
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("`${7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(1)
    $rng1.InsertBefore("`${1} `${2} `${3} `${4} `${5} `${6} ")
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: remove it from the end of the document and
#    wrap it around "${18} " instead.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("`${18} ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $d.Bookmarks.Add("_GoBack", $rng2)
}

# ---------------------------------------------------------------------------
# 3) Merge the "<w:tab/>" run with the following "Tác giả:  " run so both
#    live inside a single <w:r>, keeping the trailing hyperlink-text run and
#    the paragraph/run formatting attributes intact.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("`tTác giả:  #div_much_answer > div:nth-child(1) > div.right > div.question_info > div.author > strong > a > span", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Delete()
    $rng3.Collapse(1)
    $xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008A054D" w:rsidRDefault="008A054D" w:rsidP="007E7B2D"><w:r><w:tab/><w:t xml:space="preserve">Tác giả:  </w:t></w:r><w:r w:rsidRPr="008A054D"><w:t>#div_much_answer &gt; div:nth-child(1) &gt; div.right &gt; div.question_info &gt; div.author &gt; strong &gt; a &gt; span</w:t></w:r></w:p>'
    $rng3.InsertXML($xml3) | Out-Null
}
